$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.852.57"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "2.347.34"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("E4").Value = "  +0.11%  "
$s = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.672"
$ws.Range("D5").Style = $s
$ws.Range("E5").Value = "  -0.94%  "
$s = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "239.65"
$ws.Range("D6").Style = $s
$ws.Range("E6").Value = "  -0.90%  "
$s = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.94"
$ws.Range("D7").Style = $s
$ws.Range("E7").Value = "  -1.62%  "
$ws.Range("E8").Value = "  +0.03%  "
$s = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.593"
$ws.Range("D9").Style = $s
$ws.Range("E9").Value = "  +4.32%  "
$ws.Range("E10").Value = "  -1.68%  "
$s = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.87"
$ws.Range("D11").Style = $s
$ws.Range("E11").Value = "  +2.53%  "
$s = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "32.84"
$ws.Range("D12").Style = $s
$ws.Range("E12").Value = "  +3.39%  "
$s = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.29"
$ws.Range("D13").Style = $s
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("D15").Value = "2.697.39"
$ws.Range("E15").Value = "  -0.29%  "
$s = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.35"
$ws.Range("D16").Style = $s
$ws.Range("E16").Value = "  -2.50%  "
$s = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.904"
$ws.Range("D17").Style = $s
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("D18").Value = "2.340.95"
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("D19").Value = "43.719.95"
$ws.Range("E19").Value = "  -1.21%  "
$ws.Range("E20").Value = "  -0.90%  "
$ws.Range("E21").Value = "  +0.51%  "
$s = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "77.26"
$ws.Range("D22").Style = $s
$ws.Range("E22").Value = "  -1.08%  "
$s = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "256.26"
$ws.Range("D23").Style = $s
$ws.Range("E23").Value = "  +0.28%  "
$s = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.92"
$ws.Range("D24").Style = $s
$ws.Range("E24").Value = "  +20.60%  "
$ws.Range("E25").Value = "  -0.03%  "
$s = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.75"
$ws.Range("D26").Style = $s
$ws.Range("E26").Value = "  -0.46%  "
$s = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.51"
$ws.Range("D27").Style = $s
$ws.Range("E27").Value = "  -2.04%  "
$s = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.61"
$ws.Range("D28").Style = $s
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("E29").Value = "  -1.22%  "
$s = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.57"
$ws.Range("D30").Style = $s
$ws.Range("E30").Value = "  +0.18%  "
$s = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "177.18"
$ws.Range("D31").Style = $s
$ws.Range("E31").Value = "  +1.38%  "
$s = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.130"
$ws.Range("D32").Style = $s
$ws.Range("E32").Value = "  -0.14%  "
$s = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.136"
$ws.Range("D33").Style = $s
$ws.Range("E33").Value = "  +2.68%  "
$s = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0756"
$ws.Range("D34").Style = $s
$ws.Range("E34").Value = "  -0.55%  "
$s = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.52"
$ws.Range("D35").Style = $s
$ws.Range("E35").Value = "  +3.29%  "
$ws.Range("E36").Value = "  -3.93%  "
$s = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.82"
$ws.Range("D37").Style = $s
$ws.Range("E37").Value = "  -1.26%  "
$ws.Range("E38").Value = "  -3.20%  "
$ws.Range("E39").Value = "  -4.00%  "
$s = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0281"
$ws.Range("D40").Style = $s
$ws.Range("E40").Value = "  +2.84%  "
$s = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.31"
$ws.Range("D41").Style = $s
$ws.Range("E41").Value = "  +32.24%  "
$s = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.111"
$ws.Range("D42").Style = $s
$ws.Range("E42").Value = "  +11.20%  "
$ws.Range("E43").Value = "  +1.70%  "
$s = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.06"
$ws.Range("D44").Style = $s
$ws.Range("E44").Value = "  -1.16%  "
$s = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.202"
$ws.Range("D45").Style = $s
$ws.Range("E45").Value = "  +6.29%  "
$s = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.88"
$ws.Range("D46").Style = $s
$ws.Range("E46").Value = "  +9.04%  "
$ws.Range("E47").Value = "  +0.63%  "
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("E49").Value = "  -1.42%  "
$s = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "99.08"
$ws.Range("D50").Style = $s
$ws.Range("E50").Value = "  -1.49%  "
$ws.Range("E51").Value = "  -0.79%  "
